$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from existing header cell (H1) to the new header cells so
# they pick up the same bold/border/center style, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new I0 / IF data columns for rows 2-19.
$data = @(
    @(9,9),
    @(8,8),
    @(7,7),
    @(6,6),
    @(7,8),
    @(8,8),
    @(6,7),
    @(8,8),
    @(3,5),
    @(9,9),
    @(8,8),
    @(7,7),
    @(6,6),
    @(6,6),
    @(6,6),
    @(7,7),
    @(6,6),
    @(5,5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
